# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Swap country names (and therefore figures) for the three-way
# Jordania / Tunez / Azerbaiyan block (rows 73-75) ---
$ws.Range("A73").Value = "Azerbaiyan"
$ws.Range("A74").Value = "Jordania"
$ws.Range("A75").Value = "Tunez"

# --- Swap country names for Montserrat / Islas Malvinas (rows 216-217) ---
$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("A217").Value = "Montserrat"

# --- Update numeric figures (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Alemania (row 20)
$ws.Range("B20").Value = 417952
$ws.Range("C20").Value = 602
$ws.Range("E20").Value = 97662

# Arabia Saudita (row 25)
$ws.Range("B25").Value = 344552
$ws.Range("C25").Value = 395
$ws.Range("D25").Value = 330995
$ws.Range("E25").Value = 8276
$ws.Range("G25").Value = 17
$ws.Range("H25").Value = 5281

# Catar (row 39)
$ws.Range("B39").Value = 130965
$ws.Range("C39").Value = 254
$ws.Range("D39").Value = 127868
$ws.Range("E39").Value = 2868

# Kuwait (row 43)
$ws.Range("B43").Value = 120927
$ws.Range("C43").Value = 695
$ws.Range("D43").Value = 112110
$ws.Range("E43").Value = 8073
$ws.Range("G43").Value = 4
$ws.Range("H43").Value = 744

# Barein (row 58)
$ws.Range("E58").Value = 3119
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = 312

# Kirguistan (row 70)
$ws.Range("B70").Value = 55144
$ws.Range("C70").Value = 556
$ws.Range("D70").Value = 47337
$ws.Range("E70").Value = 6677
$ws.Range("G70").Value = 4
$ws.Range("H70").Value = 1130

# Row 73 - now Azerbaiyan
$ws.Range("B73").Value = 49013
$ws.Range("C73").Value = 792
$ws.Range("D73").Value = 41051
$ws.Range("E73").Value = 7298
$ws.Range("G73").Value = 8
$ws.Range("H73").Value = 664

# Row 74 - now Jordania
$ws.Range("B74").Value = 48930
$ws.Range("C74").Value = 0
$ws.Range("D74").Value = 7449
$ws.Range("E74").Value = 40973
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 508

# Row 75 - now Tunez
$ws.Range("B75").Value = 48799
$ws.Range("C75").Value = 1585
$ws.Range("D75").Value = 5032
$ws.Range("E75").Value = 42948
$ws.Range("G75").Value = 35
$ws.Range("H75").Value = 819

# Dinamarca (row 81)
$ws.Range("B81").Value = 39411
$ws.Range("C81").Value = 789
$ws.Range("D81").Value = 31701
$ws.Range("E81").Value = 7010
$ws.Range("G81").Value = 3
$ws.Range("H81").Value = 700

# Finlandia (row 102)
$ws.Range("B102").Value = 14652
$ws.Range("C102").Value = 178
$ws.Range("E102").Value = 4499

# Tayikistan (row 111)
$ws.Range("B111").Value = 10736
$ws.Range("C111").Value = 41
$ws.Range("D111").Value = 9836
$ws.Range("E111").Value = 819

# Vietnam (row 168)
$ws.Range("B168").Value = 1160
$ws.Range("C168").Value = 12
$ws.Range("D168").Value = 1051
$ws.Range("E168").Value = 74

# Row 216 - now Islas Malvinas
$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0

# Row 217 - now Montserrat
$ws.Range("D217").Value = 12
$ws.Range("H217").Value = 1
